# Update the AIC model-name labels in column A (rows 2-11) of the
# "combined_r" results sheet so that the abbreviated model codes are
# replaced with fully spelled-out, more interpretable model descriptions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "rugosity + year + site"
$ws.Range("A3").Value  = "year + site"
$ws.Range("A4").Value  = "year + site + year*site"
$ws.Range("A5").Value  = "site"
$ws.Range("A6").Value  = "rugosity + site"
$ws.Range("A7").Value  = "rugosity + year"
$ws.Range("A8").Value  = "rugosity + site + site*rugosity"
$ws.Range("A9").Value  = "rugosity + year + year*rugosity"
$ws.Range("A10").Value = "rugosity"
$ws.Range("A11").Value = "year"

# Column A now holds longer text, so widen it to fit (best-fit width).
$ws.Columns.Item(1).ColumnWidth = 27.833333333333332

# Reset the view back to the top-left cell (clears the old multi-cell
# selection that was left over from editing).
$ws.Range("A1").Select() | Out-Null
